# Update the column header names for all campus mapping sheets
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "building_name"
$ws.Range("B1").Value = "building_no"
$ws.Range("C1").Value = "website:map"

# The url column's displayed text had stray leading/trailing spaces around the
# sharepoi number (e.g. "...sharepoi= 898 "); trim them to "...sharepoi=898".
$poiNumbers = @("898", "901", "899", "902", "903", "905", "906", "907", "908", "911", "914", "916", "917")

for ($i = 0; $i -lt $poiNumbers.Count; $i++) {
    $row = $i + 2
    $url = "https://use.mazemap.com/?campusid=218&sharepoitype=identifier&sharepoi=" + $poiNumbers[$i]
    $ws.Cells.Item($row, 3).Value = $url
}
